$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 9438
    3 = 9438
    4 = 9438
    5 = 9438
    6 = 9438
    7 = 8764
    8 = 8764
    9 = 8764
    10 = 8764
    11 = 8764
    12 = 8764
    13 = 8764
    14 = 8764
    15 = 8764
    16 = 8764
    17 = 8764
    18 = 8702
    19 = 8702
    20 = 8702
    21 = 8702
    22 = 8702
    23 = 8257
    24 = 8257
    25 = 8257
    26 = 8257
    27 = 8257
    28 = 8257
    29 = 8257
    30 = 8257
    31 = 8257
    32 = 8257
    33 = 8257
    34 = 8173
    35 = 8173
    36 = 8173
    37 = 8173
    38 = 8173
    39 = 8173
    40 = 8173
    41 = 8173
    42 = 8173
    43 = 8173
    44 = 8173
    45 = 8165
    46 = 8165
    47 = 8165
    48 = 8165
    49 = 8165
    50 = 8165
    51 = 8165
    52 = 8165
    53 = 8165
    54 = 8165
    55 = 8165
    56 = 8165
    57 = 8165
    58 = 8165
    59 = 8165
    60 = 8165
    61 = 8165
    62 = 8165
    63 = 8165
    64 = 7586
    65 = 7586
    66 = 7586
    67 = 7586
    68 = 7586
    69 = 7586
    70 = 7586
    71 = 7586
    72 = 7586
    73 = 7586
    74 = 7586
    75 = 7586
    76 = 7586
    77 = 7586
    78 = 7586
    79 = 7586
    80 = 7586
    81 = 7586
    82 = 7586
    83 = 7586
    84 = 7586
    85 = 7586
    86 = 7586
    87 = 7586
    88 = 7586
    89 = 7586
    90 = 7586
    91 = 7586
    92 = 7586
    93 = 7586
    94 = 7586
    95 = 7586
    96 = 7586
    97 = 7586
    98 = 7586
    99 = 7586
    100 = 7586
    101 = 7586
    102 = 7586
    103 = 7586
    104 = 7586
    105 = 7586
    106 = 7586
    107 = 7586
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
